$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after the header (rows 2-4), shifting existing data down
$ws.Range("A2:C4").Insert()
$ws.Range("A2:C4").ClearFormats()

$topData = @(
    @(-9.328194618225098, -6.853240489959717, 5.245095252990723),
    @(-7.405017375946045, -14.76152610778809, 0.4900901317596435),
    @(-8.944772720336914, -6.256664276123047, 6.951004505157471)
)
for ($i = 0; $i -lt $topData.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topData[$i][0]
    $ws.Cells.Item($r, 2).Value = $topData[$i][1]
    $ws.Cells.Item($r, 3).Value = $topData[$i][2]
}

# Append 7 new rows of data at the bottom (rows 25-31)
$bottomData = @(
    @(4.273637771606445, -40.02112579345703, 19.01494979858398),
    @(-33.29051208496094, 13.78821468353272, -5.866414546966553),
    @(-38.59527587890625, -11.46046161651611, -14.89169502258301),
    @(69.93356323242188, -70.74887084960938, 44.17119979858398),
    @(-14.73312759399414, 12.10027313232422, 5.143325805664063),
    @(16.71687316894531, -16.55658340454102, 22.30324172973633),
    @(-58.74892807006836, -18.46404838562012, -22.14792823791504)
)
for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomData[$i][2]
}
